$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.342.13'
$ws.Range('E2').Value = '  -1.36%  '

$ws.Range('D3').Value = '2.051.57'
$ws.Range('E3').Value = '  -1.54%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.14'
$ws.Range('E5').Value = '  -0.95%  '

$ws.Range('E6').Value = '  -0.66%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.00'
$ws.Range('E8').Value = '  -3.89%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0769'
$ws.Range('E10').Value = '  -2.59%  '

$ws.Range('E11').Value = '  +1.30%  '

$ws.Range('D12').Value = '2.356.16'
$ws.Range('E12').Value = '  -1.46%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.61'
$ws.Range('E13').Value = '  -1.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.60'
$ws.Range('E14').Value = '  -3.03%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.756'
$ws.Range('E15').Value = '  -2.69%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.26'
$ws.Range('E16').Value = '  -1.81%  '

$ws.Range('D17').Value = '2.040.70'
$ws.Range('E17').Value = '  -1.45%  '

$ws.Range('D18').Value = '37.280.89'
$ws.Range('E18').Value = '  -1.31%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.06'
$ws.Range('E19').Value = '  -1.86%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.66'
$ws.Range('E20').Value = '  -2.64%  '

$ws.Range('D21').Value = '0.0₃0823'
$ws.Range('E21').Value = '  -3.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.52'
$ws.Range('E22').Value = '  -0.87%  '

$ws.Range('E23').Value = '  +0.09%  '

$ws.Range('E24').Value = '  -0.11%  '

$ws.Range('E25').Value = '  -3.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.81'
$ws.Range('E26').Value = '  +6.87%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.86'
$ws.Range('E27').Value = '  -1.15%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.129'
$ws.Range('E28').Value = '  -5.28%  '

$ws.Range('E29').Value = '  -1.70%  '

$ws.Range('E30').Value = '  -5.83%  '

$ws.Range('E31').Value = '  -0.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.52'
$ws.Range('E32').Value = '  -4.47%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0621'
$ws.Range('E33').Value = '  -1.99%  '

$ws.Range('E34').Value = '  -4.09%  '

$ws.Range('E35').Value = '  -2.26%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.27'
$ws.Range('E37').Value = '  -4.88%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.20%  '

$ws.Range('E39').Value = '  -2.07%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0225'
$ws.Range('E40').Value = '  +3.33%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.08'
$ws.Range('E41').Value = '  -1.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0952'
$ws.Range('E42').Value = '  -3.42%  '

$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.90'
$ws.Range('E43').Value = '  +0.45%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.478.05'
$ws.Range('E44').Value = '  +2.19%  '

$ws.Range('E45').Value = '  +1.55%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.56'
$ws.Range('E46').Value = '  -0.47%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.99'
$ws.Range('E47').Value = '  -4.57%  '

$ws.Range('E48').Value = '  -3.26%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.24'
$ws.Range('E49').Value = '  -2.24%  '

$ws.Range('E50').Value = '  -2.10%  '

$ws.Range('D51').Value = '2.242.04'
$ws.Range('E51').Value = '  -1.48%  '
